$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value = 3400
$ws.Range("I62").Value = 1500
$ws.Range("J62").Value = 3875
$ws.Range("K62").Value = 1500
$ws.Range("L62").Value = 3875
$ws.Range("M62").Value = -876
$ws.Range("N62").Value = -5123
# row 65
$ws.Range("H65").Value = 3400
$ws.Range("I65").Value = 1500
$ws.Range("J65").Value = 3875
$ws.Range("K65").Value = 7500
$ws.Range("L65").Value = 19375
$ws.Range("M65").Value = -4380
$ws.Range("N65").Value = -25615
# row 99
$ws.Range("H99").Value = 1266.7333
$ws.Range("I99").Value = 357.5
$ws.Range("J99").Value = 2305.8572
$ws.Range("K99").Value = 1072.5
$ws.Range("L99").Value = 6917.571599999999
$ws.Range("M99").Value = 425.5
$ws.Range("N99").Value = -9913.571599999999
# row 100
$ws.Range("H100").Value = 2034
$ws.Range("J100").Value = 2176.5
$ws.Range("L100").Value = 2176.5
$ws.Range("N100").Value = -3258.5
# row 113
$ws.Range("H113").Value = 3207.7827
$ws.Range("I113").Value = 2998.75
$ws.Range("J113").Value = 3685.5715
$ws.Range("K113").Value = 2998.75
$ws.Range("L113").Value = 3685.5715
$ws.Range("M113").Value = 255.25
$ws.Range("N113").Value = -10193.5715
# row 129
$ws.Range("H129").Value = 821.625
$ws.Range("J129").Value = 1158.4286
$ws.Range("L129").Value = 3475.2858
$ws.Range("N129").Value = -13475.2858
# row 132
$ws.Range("H132").Value = 629772.5
$ws.Range("I132").Value = 1377.8923
$ws.Range("J132").Value = 3771745.5
$ws.Range("K132").Value = 4133.6769
$ws.Range("L132").Value = 11315236.5
$ws.Range("M132").Value = -1603.6769
$ws.Range("N132").Value = -11320296.5
# row 138
$ws.Range("H138").Value = 2166426
$ws.Range("I138").Value = 1382.6666
$ws.Range("K138").Value = 4147.9998
$ws.Range("M138").Value = 992.0002000000004

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 43566440
$ws.Range("I61").Value = 50051330
$ws.Range("J61").Value = 333833.34
$ws.Range("K61").Value = 50051330
$ws.Range("L61").Value = 333833.34
$ws.Range("M61").Value = -50051118
$ws.Range("N61").Value = -334257.34
# row 74
$ws.Range("H74").Value = 7413503
$ws.Range("I74").Value = 11954064
$ws.Range("J74").Value = 78749.84
$ws.Range("K74").Value = 11954064
$ws.Range("L74").Value = 78749.84
$ws.Range("M74").Value = -11953190
$ws.Range("N74").Value = -80497.84
# row 77
$ws.Range("H77").Value = 7413503
$ws.Range("I77").Value = 11954064
$ws.Range("J77").Value = 78749.84
$ws.Range("K77").Value = 59770320
$ws.Range("L77").Value = 393749.2
$ws.Range("M77").Value = -59765952
$ws.Range("N77").Value = -402485.2
# row 97
$ws.Range("H97").Value = 2404661.2
$ws.Range("I97").Value = 3472982.8
$ws.Range("J97").Value = 937.5
$ws.Range("K97").Value = 3472982.8
$ws.Range("L97").Value = 937.5
$ws.Range("M97").Value = -3472486.8
$ws.Range("N97").Value = -1929.5
# row 102
$ws.Range("H102").Value = 7149573
$ws.Range("I102").Value = 7943692
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 7943692
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -7942070
$ws.Range("N102").Value = -5744
# row 132
$ws.Range("H132").Value = 262374.5
$ws.Range("I132").Value = 201199.2
$ws.Range("J132").Value = 364333.34
$ws.Range("K132").Value = 603597.6000000001
$ws.Range("L132").Value = 1093000.02
$ws.Range("M132").Value = -601067.6000000001
$ws.Range("N132").Value = -1098060.02
# row 136
$ws.Range("H136").Value = 43566440
$ws.Range("I136").Value = 50051330
$ws.Range("J136").Value = 333833.34
$ws.Range("K136").Value = 150153990
$ws.Range("L136").Value = 1001500.02
$ws.Range("M136").Value = -150151440
$ws.Range("N136").Value = -1006600.02

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 556.5833
$ws.Range("I94").Value = 334.875
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 334.875
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 116.125
$ws.Range("N94").Value = -1902
# row 99
$ws.Range("H99").Value = 1240.3334
$ws.Range("I99").Value = 1105
$ws.Range("K99").Value = 1105
$ws.Range("M99").Value = 393
# row 107
$ws.Range("H107").Value = 1583.1
$ws.Range("I107").Value = 1135.2667
$ws.Range("K107").Value = 1135.2667
$ws.Range("M107").Value = 784.7333000000001
# row 134
$ws.Range("H134").Value = 2123.543
$ws.Range("I134").Value = 1060.68
$ws.Range("J134").Value = 4780.7
$ws.Range("K134").Value = 3182.04
$ws.Range("L134").Value = 14342.1
$ws.Range("M134").Value = -647.04
$ws.Range("N134").Value = -19412.1
# row 141
$ws.Range("H141").Value = 45745.43
$ws.Range("J141").Value = 47760
$ws.Range("L141").Value = 47760
$ws.Range("N141").Value = -58120

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
# row 140
$ws.Range("H140").Value = 2337.541
$ws.Range("I140").Value = 2382.3809
$ws.Range("J140").Value = 2314
$ws.Range("K140").Value = 7147.1427
$ws.Range("L140").Value = 6942
$ws.Range("M140").Value = -1967.1427
$ws.Range("N140").Value = -17302

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 183795
$ws.Range("I132").Value = 144827.86
$ws.Range("J132").Value = 251987.5
$ws.Range("K132").Value = 434483.58
$ws.Range("L132").Value = 755962.5
$ws.Range("M132").Value = -431953.58
$ws.Range("N132").Value = -761022.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 740.3226
$ws.Range("I46").Value = 674.2857
$ws.Range("J46").Value = 794.7059
$ws.Range("K46").Value = 674.2857
$ws.Range("L46").Value = 794.7059
$ws.Range("M46").Value = -486.2857
$ws.Range("N46").Value = -1170.7059
# row 93
$ws.Range("H93").Value = 963.1177
$ws.Range("I93").Value = 918.5454999999999
$ws.Range("J93").Value = 1044.8334
$ws.Range("K93").Value = 918.5454999999999
$ws.Range("L93").Value = 1044.8334
$ws.Range("M93").Value = 329.4545000000001
$ws.Range("N93").Value = -3540.8334
# row 132
$ws.Range("H132").Value = 27111.05
$ws.Range("I132").Value = 17666.492
$ws.Range("K132").Value = 52999.476
$ws.Range("M132").Value = -50469.476

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 109
$ws.Range("H109").Value = 25000
$ws.Range("J109").Value = 25000
$ws.Range("L109").Value = 25000
$ws.Range("N109").Value = -27774
